# Quest.xlsx — add the "invite reward" quests (4001-4004) to Sheet1,
# bump the height of the existing data rows, and move the selection
# cursor down to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-15 (the existing quest rows) pick up the same row height used
# elsewhere in the sheet.
for ($r = 2; $r -le 15; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.85
}

# New quest definitions: invite-reward checkpoints (id 4001-4004).
$newRows = @(
    @{ Row = 16; Id = 4001; Reward = "[[1,100]]"; Note = "被邀请人奖励" },
    @{ Row = 17; Id = 4002; Reward = "[[1,80]]";  Note = "邀请人奖励" },
    @{ Row = 18; Id = 4003; Reward = "[[1,8]]";   Note = "邀请人父节点奖励" },
    @{ Row = 19; Id = 4004; Reward = "[[1,82]]";  Note = "邀请人祖父节点奖励" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Id   # id
    $ws.Cells.Item($r, 2).Value = 1          # checkpoint
    $ws.Cells.Item($r, 3).Value = 4          # type
    $ws.Cells.Item($r, 4).Value = $item.Reward  # reward
    $ws.Cells.Item($r, 6).Value = $item.Note    # 备注
    $ws.Cells.Item($r, 7).Value = 0             # settlement_type
}

# Leave the cursor where the author left off.
$ws.Range("F21").Select()
